# MasterServer.xlsx - "Property" sheet: fill in the first data row (row 2)
# with the concrete MasterServer_1 entry, and move the active selection to H3.
#
# Column layout (row 1 headers): A=ID, B=ServerID, C=Name, D=MaxOnline,
# E=CpuCount, F=IP, G=Port.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Write the text columns in the same order Excel would assign new shared
# strings (ServerID, then IP, then the ID/Name pair) so the shared-string
# table comes out in the expected order.
$ws.Range("B2").Value = "000106001"
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "MasterServer_1"

# C2 is a brand-new cell; give it the same "text" number format used by the
# other text cells in this row (A2/B2) before writing its value.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "MasterServer_1"

# Numeric columns.
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 2001

# Move the selection to H3 (was A2:H7 with active cell H7).
$ws.Range("H3").Select()
